$d = $word.ActiveDocument

# Mark both inline pictures as "no proof" (matches the rPr/noProof Word
# stamps on picture runs when it re-saves the document).
foreach ($shp in $d.InlineShapes) {
    $shp.Range.NoProofing = 1
}

# Replace the empty paragraph after "Github repo link:" with one containing
# a hyperlink to the GitHub repo.
$last = $d.Paragraphs.Last
$d.Hyperlinks.Add($last.Range, "https://github.com/AHenryHub/CSC-MS-Assignment2/tree/assignment2-henry") | Out-Null
